# "docs(readme): update timeline and gantt"
#
# Updates the Gantt chart on the "Project Planner" sheet:
#  - tweaks a few Actual-Start/Actual-Duration/%-Complete figures for
#    already-listed 3.x implementation tasks
#  - inserts four new implementation tasks (3.18-3.21) right after
#    "3.17 Tutorial pages", pushing every following row down by 4
#  - nudges the %-complete figure for "4.2 Final Report"
#  - re-anchors the conditional formatting / view state to the new extent

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# --- row 16 ("3. Implementation" section header): record an actual-start ---
$ws.Range("E16").Value = 14
$ws.Range("E16").HorizontalAlignment = -4108   # xlCenter, matches sibling C16/D16

# --- existing 3.x rows: progress updates ---
$ws.Range("F28").Value = 4        # 3.12 animate merge and heap - Actual Duration
$ws.Range("G28").Value = 1        # 3.12 animate merge and heap - % Complete
$ws.Range("G31").Value = 0.9      # 3.15 correctness tutorial  - % Complete
$ws.Range("G33").Value = 0.66     # 3.17 Tutorial pages        - % Complete

# --- insert four new task rows (3.18 - 3.21) right after row 33 ---
$ws.Range("B34:B37").EntireRow.Insert()
$ws.Range("B34:B37").RowHeight = 30

$ws.Range("B34").Value = "3.18 Procedure imple"
$ws.Range("C34").Value = 23
$ws.Range("D34").Value = 2
$ws.Range("E34").Value = 23
$ws.Range("G34").Value = 0.6

$ws.Range("B35").Value = "3.19 Improve transition anime"
$ws.Range("C35").Value = 24
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 24

$ws.Range("B36").Value = "3.20 Refactor previous work"
$ws.Range("C36").Value = 24
$ws.Range("D36").Value = 1
$ws.Range("E36").Value = 24

$ws.Range("B37").Value = "3.21 Correctness proof"
$ws.Range("C37").Value = 24
$ws.Range("D37").Value = 1

# --- "4.2 Final Report" (now row 44 after the insert above): % Complete ---
$ws.Range("G44").Value = 0.25

# --- re-anchor conditional formatting to the grown data range ---
$cfs = $ws.Cells.FormatConditions()
$cfs.Item(1).ModifyAppliesToRange($ws.Range("H5:AI53"))
$cfs.Item(9).ModifyAppliesToRange($ws.Range("B54:BO54"))

# --- view state: selection / zoom match the author's saved snapshot ---
$ws.Range("F36").Select() | Out-Null
$excel.ActiveWindow.Zoom = 89
